$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to reflect the
# repulled / recalculated data as described in the commit message.
$ws.Range("F2").Value  = -5
$ws.Range("F3").Value  = -5
$ws.Range("F5").Value  = -4
$ws.Range("F9").Value  = -3
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -3
$ws.Range("F15").Value = -11
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = 5
$ws.Range("F24").Value = -5
